$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.978.30'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '3.505.40'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'577.49"
$ws.Range("E5").Value = '  +5.04%  '
$ws.Range("D6").Value = "'178.32"
$ws.Range("E6").Value = '  -5.49%  '
$ws.Range("D7").Value = "'0.635"
$ws.Range("E7").Value = '  +4.70%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = "'0.634"
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("E10").Value = '  +3.86%  '
$ws.Range("D11").Value = "'55.45"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = "'0.0000273"
$ws.Range("E12").Value = '  +2.00%  '
$ws.Range("D13").Value = "'9.24"
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").Value = '4.069.61'
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("D15").Value = '3.504.62'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = "'18.37"
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = "'12.08"
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '65.910.96'
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").Value = "'414.46"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = '  +7.80%  '
$ws.Range("D23").Value = "'4.32"
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").Value = "'85.86"
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("E25").Value = '  +11.83%  '
$ws.Range("D26").Value = "'11.01"
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").Value = "'2.86"
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").Value = "'9.11"
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("D29").Value = "'30.45"
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").Value = "'627.24"
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("D32").Value = "'11.68"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("E34").Value = '  +13.94%  '
$ws.Range("D35").Value = "'59.63"
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '0.0₃0798'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("D38").Value = "'37.18"
$ws.Range("E38").Value = '  -3.86%  '
$ws.Range("D39").Value = "'3.51"
$ws.Range("E39").Value = '  +4.44%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = "'0.381"
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.275.26'
$ws.Range("E41").Value = '  +8.84%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = "'2.92"
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").Value = "'3.28"
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").Value = "'0.0418"
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  -4.97%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").Value = "'0.133"
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = "'8.61"
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'140.21"
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("D51").Value = "'2.38"
$ws.Range("E51").Value = '  -1.39%  '
